$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H33").Value = 243.16667
$ws.Range("J33").Value = 203.125
$ws.Range("L33").Value = 203.125
$ws.Range("N33").Value = -661.125
$ws.Range("H40").Value = 6431.619
$ws.Range("I40").Value = 4395.1714
$ws.Range("K40").Value = 4395.1714
$ws.Range("M40").Value = -4220.1714
$ws.Range("H64").Value = 3601
$ws.Range("J64").Value = 3901.5
$ws.Range("L64").Value = 3901.5
$ws.Range("N64").Value = -4397.5
$ws.Range("H67").Value = 3601
$ws.Range("J67").Value = 3901.5
$ws.Range("L67").Value = 3901.5
$ws.Range("N67").Value = -5617.5
$ws.Range("H74").Value = 9499
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 9499
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 9499
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -11371
$ws.Range("H77").Value = 9499
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 9499
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 47495
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -56855
$ws.Range("H80").Value = 13159325
$ws.Range("I80").Value = 25001106
$ws.Range("J80").Value = 1789.1111
$ws.Range("K80").Value = 75003318
$ws.Range("L80").Value = 5367.3333
$ws.Range("M80").Value = -75002320
$ws.Range("N80").Value = -7363.3333
$ws.Range("H83").Value = 13159325
$ws.Range("I83").Value = 25001106
$ws.Range("J83").Value = 1789.1111
$ws.Range("K83").Value = 225009954
$ws.Range("L83").Value = 16101.9999
$ws.Range("M83").Value = -225004962
$ws.Range("N83").Value = -26085.9999
$ws.Range("H86").Value = 45456440
$ws.Range("I86").Value = 50001964
$ws.Range("K86").Value = 50001964
$ws.Range("M86").Value = -50000841
$ws.Range("H89").Value = 45456440
$ws.Range("I89").Value = 50001964
$ws.Range("K89").Value = 250009820
$ws.Range("M89").Value = -250004204
$ws.Range("H107").Value = 951.6389
$ws.Range("I107").Value = 944.875
$ws.Range("K107").Value = 944.875
$ws.Range("M107").Value = 975.125
$ws.Range("H112").Value = 3708.5881
$ws.Range("J112").Value = 3767.0303
$ws.Range("L112").Value = 11301.0909
$ws.Range("N112").Value = -13517.0909
$ws.Range("H131").Value = 5050.769
$ws.Range("I131").Value = 1542.091
$ws.Range("K131").Value = 4626.272999999999
$ws.Range("M131").Value = 413.7270000000008
$ws.Range("H132").Value = 191239.61
$ws.Range("I132").Value = 208424.95
$ws.Range("J132").Value = 6497.25
$ws.Range("K132").Value = 625274.8500000001
$ws.Range("L132").Value = 19491.75
$ws.Range("M132").Value = -622744.8500000001
$ws.Range("N132").Value = -24551.75
$ws.Range("H135").Value = 5464.615
$ws.Range("I135").Value = 935.5714
$ws.Range("K135").Value = 8420.142600000001
$ws.Range("M135").Value = -5885.142600000001
$ws.Range("H137").Value = 3132.963
$ws.Range("I137").Value = 1998.8
$ws.Range("K137").Value = 5996.4
$ws.Range("M137").Value = -3446.4
$ws.Range("H138").Value = 5164.6045
$ws.Range("I138").Value = 2704.611
$ws.Range("J138").Value = 5815.7793
$ws.Range("K138").Value = 8113.833
$ws.Range("L138").Value = 17447.3379
$ws.Range("M138").Value = -2973.833
$ws.Range("N138").Value = -27727.3379
$ws.Range("H139").Value = 135840.47
$ws.Range("J139").Value = 135840.47
$ws.Range("L139").Value = 135840.47
$ws.Range("N139").Value = -146120.47
$ws.Range("H140").Value = 99157.2
$ws.Range("J140").Value = 99157.2
$ws.Range("L140").Value = 99157.2
$ws.Range("N140").Value = -109517.2
$ws.Range("H141").Value = 2589
$ws.Range("J141").Value = 4104.5
$ws.Range("L141").Value = 12313.5
$ws.Range("N141").Value = -22673.5
# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 15062.543
$ws.Range("I32").Value = 14387.906
$ws.Range("K32").Value = 14387.906
$ws.Range("M32").Value = -14100.906
$ws.Range("H45").Value = 25900.334
$ws.Range("I45").Value = 30490.62
$ws.Range("J45").Value = 4479
$ws.Range("K45").Value = 30490.62
$ws.Range("L45").Value = 4479
$ws.Range("M45").Value = -30113.62
$ws.Range("N45").Value = -5233
$ws.Range("H61").Value = 5836.7144
$ws.Range("I61").Value = 5143.9165
$ws.Range("K61").Value = 5143.9165
$ws.Range("M61").Value = -4931.9165
$ws.Range("H74").Value = 8930616
$ws.Range("I74").Value = 10871280
$ws.Range("J74").Value = 3560
$ws.Range("K74").Value = 10871280
$ws.Range("L74").Value = 3560
$ws.Range("M74").Value = -10870406
$ws.Range("N74").Value = -5308
$ws.Range("H77").Value = 8930616
$ws.Range("I77").Value = 10871280
$ws.Range("J77").Value = 3560
$ws.Range("K77").Value = 54356400
$ws.Range("L77").Value = 17800
$ws.Range("M77").Value = -54352032
$ws.Range("N77").Value = -26536
$ws.Range("H108").Value = 80000
$ws.Range("J108").Value = 80000
$ws.Range("L108").Value = 80000
$ws.Range("N108").Value = -87680
$ws.Range("H110").Value = 1378.5526
$ws.Range("I110").Value = 1371.8055
$ws.Range("K110").Value = 1371.8055
$ws.Range("M110").Value = 673.1945000000001
$ws.Range("H132").Value = 15564.236
$ws.Range("I132").Value = 18943.584
$ws.Range("K132").Value = 56830.75199999999
$ws.Range("M132").Value = -54300.75199999999
$ws.Range("H136").Value = 5836.7144
$ws.Range("I136").Value = 5143.9165
$ws.Range("K136").Value = 15431.7495
$ws.Range("M136").Value = -12881.7495
# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H81").Value = 51282.832
$ws.Range("J81").Value = 59397.6
$ws.Range("L81").Value = 59397.6
$ws.Range("N81").Value = -61519.6
$ws.Range("H84").Value = 51282.832
$ws.Range("J84").Value = 59397.6
$ws.Range("L84").Value = 178192.8
$ws.Range("N84").Value = -188800.8
$ws.Range("H105").Value = 2618.0625
$ws.Range("I105").Value = 1521.3
$ws.Range("K105").Value = 1521.3
$ws.Range("M105").Value = 225.7
$ws.Range("H132").Value = 117696
$ws.Range("J132").Value = 117696
$ws.Range("L132").Value = 117696
$ws.Range("N132").Value = -127816
$ws.Range("H134").Value = 4788.222
$ws.Range("I134").Value = 2892.125
$ws.Range("K134").Value = 8676.375
$ws.Range("M134").Value = -6141.375
# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value = 104.958336
$ws.Range("I7").Value = 61.909092
$ws.Range("J7").Value = 141.38461
$ws.Range("K7").Value = 61.909092
$ws.Range("L7").Value = 141.38461
$ws.Range("M7").Value = 51.090908
$ws.Range("N7").Value = -367.38461
$ws.Range("H22").Value = 462.58334
$ws.Range("I22").Value = 341.66666
$ws.Range("J22").Value = 583.5
$ws.Range("K22").Value = 341.66666
$ws.Range("L22").Value = 583.5
$ws.Range("M22").Value = 8.333340000000021
$ws.Range("N22").Value = -1283.5
$ws.Range("H31").Value = 14708757
$ws.Range("I31").Value = 26318046
$ws.Range("J31").Value = 3656.5334
$ws.Range("K31").Value = 26318046
$ws.Range("L31").Value = 3656.5334
$ws.Range("M31").Value = -26317751
$ws.Range("N31").Value = -4246.5334
$ws.Range("H34").Value = 14708757
$ws.Range("I34").Value = 26318046
$ws.Range("J34").Value = 3656.5334
$ws.Range("K34").Value = 26318046
$ws.Range("L34").Value = 3656.5334
$ws.Range("M34").Value = -26317844
$ws.Range("N34").Value = -4060.5334
$ws.Range("H58").Value = 2260.875
$ws.Range("I58").Value = 2098.1428
$ws.Range("K58").Value = 2098.1428
$ws.Range("M58").Value = -1895.1428
$ws.Range("H93").Value = 52652460
$ws.Range("I93").Value = 12341.714
$ws.Range("J93").Value = 200044800
$ws.Range("K93").Value = 12341.714
$ws.Range("L93").Value = 200044800
$ws.Range("M93").Value = -10469.714
$ws.Range("N93").Value = -200048544
$ws.Range("H105").Value = 1567.8
$ws.Range("I105").Value = 1445.5
$ws.Range("K105").Value = 1445.5
$ws.Range("M105").Value = 301.5
$ws.Range("H107").Value = 933.56525
$ws.Range("I107").Value = 774.58826
$ws.Range("K107").Value = 774.58826
$ws.Range("M107").Value = 1145.41174
$ws.Range("H122").Value = 4037.682
$ws.Range("I122").Value = 2309.2727
$ws.Range("K122").Value = 6927.8181
$ws.Range("M122").Value = -4477.8181
$ws.Range("H134").Value = 1502.8276
$ws.Range("I134").Value = 1502.8276
$ws.Range("K134").Value = 4508.4828
$ws.Range("M134").Value = -1973.4828
$ws.Range("H136").Value = 2260.875
$ws.Range("I136").Value = 2098.1428
$ws.Range("K136").Value = 6294.428400000001
$ws.Range("M136").Value = -3744.428400000001
$ws.Range("H141").Value = 275873
$ws.Range("J141").Value = 295312.34
$ws.Range("L141").Value = 295312.34
$ws.Range("N141").Value = -305672.34
# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H34").Value = 1341.0667
$ws.Range("J34").Value = 3339
$ws.Range("L34").Value = 10017
$ws.Range("N34").Value = -10185
$ws.Range("H55").Value = 3749.9
$ws.Range("J55").Value = 4500
$ws.Range("L55").Value = 13500
$ws.Range("N55").Value = -13854
$ws.Range("H113").Value = 2760.6667
$ws.Range("J113").Value = 2828.2693
$ws.Range("L113").Value = 8484.8079
$ws.Range("N113").Value = -12824.8079
$ws.Range("H121").Value = 858
$ws.Range("I121").Value = 866.3333
$ws.Range("J121").Value = 833
$ws.Range("K121").Value = 2598.9999
$ws.Range("L121").Value = 2499
$ws.Range("M121").Value = -1288.9999
$ws.Range("N121").Value = -5119
$ws.Range("H132").Value = 2499.6428
$ws.Range("J132").Value = 3499
$ws.Range("L132").Value = 31491
$ws.Range("N132").Value = -36551
# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H53").Value = 49995.332
$ws.Range("J53").Value = 49995.332
$ws.Range("L53").Value = 49995.332
$ws.Range("N53").Value = -51257.332
$ws.Range("H63").Value = 83371.336
$ws.Range("J63").Value = 83371.336
$ws.Range("L63").Value = 83371.336
$ws.Range("N63").Value = -84743.336
$ws.Range("H66").Value = 83371.336
$ws.Range("J66").Value = 83371.336
$ws.Range("L66").Value = 250114.008
$ws.Range("N66").Value = -256978.008
$ws.Range("H102").Value = 26851196
$ws.Range("I102").Value = 42509690
$ws.Range("J102").Value = 8067
$ws.Range("K102").Value = 42509690
$ws.Range("L102").Value = 8067
$ws.Range("M102").Value = -42508068
$ws.Range("N102").Value = -11311
$ws.Range("H113").Value = 1692.1666
$ws.Range("I113").Value = 1640
$ws.Range("J113").Value = 1744.3334
$ws.Range("K113").Value = 1640
$ws.Range("L113").Value = 1744.3334
$ws.Range("M113").Value = 530
$ws.Range("N113").Value = -6084.3334
$ws.Range("H122").Value = 4071.625
$ws.Range("I122").Value = 2064
$ws.Range("J122").Value = 5714.227
$ws.Range("K122").Value = 6192
$ws.Range("L122").Value = 17142.681
$ws.Range("M122").Value = -3742
$ws.Range("N122").Value = -22042.681
$ws.Range("H126").Value = 4626.476
$ws.Range("I126").Value = 3221.4138
$ws.Range("J126").Value = 7760.846
$ws.Range("K126").Value = 9664.241399999999
$ws.Range("L126").Value = 23282.538
$ws.Range("M126").Value = -7194.241399999999
$ws.Range("N126").Value = -28222.538
$ws.Range("H132").Value = 4201.3823
$ws.Range("I132").Value = 3883.0833
$ws.Range("J132").Value = 4965.3
$ws.Range("K132").Value = 11649.2499
$ws.Range("L132").Value = 14895.9
$ws.Range("M132").Value = -9119.249899999999
$ws.Range("N132").Value = -19955.9
# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 5281.88
$ws.Range("I16").Value = 4864.143
$ws.Range("J16").Value = 7475
$ws.Range("K16").Value = 4864.143
$ws.Range("L16").Value = 7475
$ws.Range("M16").Value = -4694.143
$ws.Range("N16").Value = -7815
$ws.Range("H22").Value = 1144.9333
$ws.Range("I22").Value = 1029.9
$ws.Range("K22").Value = 1029.9
$ws.Range("M22").Value = -734.9000000000001
$ws.Range("H27").Value = 1144.9333
$ws.Range("I27").Value = 1029.9
$ws.Range("K27").Value = 1029.9
$ws.Range("M27").Value = -922.9000000000001
$ws.Range("H40").Value = 27783472
$ws.Range("I40").Value = 16671294
$ws.Range("K40").Value = 16671294
$ws.Range("M40").Value = -16671158
$ws.Range("H93").Value = 2662
$ws.Range("I93").Value = 2188
$ws.Range("J93").Value = 2899
$ws.Range("K93").Value = 2188
$ws.Range("L93").Value = 2899
$ws.Range("M93").Value = -940
$ws.Range("N93").Value = -5395
$ws.Range("H133").Value = 64953
$ws.Range("J133").Value = 64953
$ws.Range("L133").Value = 64953
$ws.Range("N133").Value = -70013
$ws.Range("H136").Value = 3996.7407
$ws.Range("I136").Value = 3540.9524
$ws.Range("K136").Value = 10622.8572
$ws.Range("M136").Value = -8072.8572
# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H37").Value = 36756
$ws.Range("I37").Value = 20026
$ws.Range("J37").Value = 42332.668
$ws.Range("K37").Value = 20026
$ws.Range("L37").Value = 42332.668
$ws.Range("M37").Value = -19823
$ws.Range("N37").Value = -42738.668
$ws.Range("H122").Value = 3713.3157
$ws.Range("I122").Value = 3124.2083
$ws.Range("J122").Value = 4723.2144
$ws.Range("K122").Value = 9372.624899999999
$ws.Range("L122").Value = 14169.6432
$ws.Range("M122").Value = -6922.624899999999
$ws.Range("N122").Value = -19069.6432
$ws.Range("H123").Value = 49989
$ws.Range("J123").Value = 49989
$ws.Range("L123").Value = 49989
$ws.Range("N123").Value = -59789
$ws.Range("H132").Value = 7410795.5
$ws.Range("I132").Value = 23813238
$ws.Range("J132").Value = 3240.9355
$ws.Range("K132").Value = 71439714
$ws.Range("L132").Value = 9722.806500000001
$ws.Range("M132").Value = -71437184
$ws.Range("N132").Value = -14782.8065
$ws.Range("H133").Value = 105000
$ws.Range("J133").Value = 105000
$ws.Range("L133").Value = 105000
$ws.Range("N133").Value = -115120
$ws.Range("H136").Value = 4391.1016
$ws.Range("I136").Value = 2867.9058
$ws.Range("J136").Value = 9436.6875
$ws.Range("K136").Value = 8603.7174
$ws.Range("L136").Value = 28310.0625
$ws.Range("M136").Value = -6053.7174
$ws.Range("N136").Value = -33410.0625

Write-Host "Applied all updates"